$d = $word.ActiveDocument

# "LOB1053: Física III (Requisito fraco)" is followed by three paragraphs
# that the site regenerate no longer emits:
#   - an empty paragraph
#   - "Ver no Jupiter Salvar em pdf Salvar em docx"
#   - "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
#      pages. Original theme under Creative Commons Attribution"
# Find the anchor paragraph, then delete the run of three paragraphs
# (including their paragraph marks) that follow it, leaving the trailing
# empty paragraph and the page-break paragraph untouched.

$anchorPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*LOB1053: Física III (Requisito fraco)*") {
        $anchorPara = $p
        break
    }
}

$firstToDelete = $anchorPara.Next()
$lastToDelete = $firstToDelete.Next().Next()

$deleteRange = $d.Range($firstToDelete.Range.Start, $lastToDelete.Range.End)
$deleteRange.Delete()
